$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.351.51'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.994.74'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.33%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.76'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.33'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +12.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +4.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.993.97'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.45%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.89'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.99%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +5.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000227'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +7.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.39'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.488.32'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.08'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +9.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.995.26'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '58.267.67'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '425.17'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.64'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +6.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.714'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +9.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.52'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +7.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.12'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +5.75%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.63'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.89%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.53'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.08'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +9.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.64'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +6.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.87'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.08'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0980'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.78'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +8.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.967'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +7.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.10'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0724'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +17.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.92'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +6.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '48.56'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +17.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '399.73'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +11.34%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.67%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.734.61'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.81%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.246'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +7.93%  '
$ws.Range('B46').NumberFormat = '@'
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').NumberFormat = '@'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('B47').NumberFormat = '@'
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').NumberFormat = '@'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '125.24'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.60%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.96%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.02%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.37'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +4.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.03'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.76%  '
